$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# NOTE on ordering: the underlying shared-strings table is append/compact
# (a string is dropped once its last reference is rewritten, and freshly
# seen strings are appended in first-seen order). The writes below are
# deliberately ordered to reproduce the exact target shared-string table:
#   45 /Planeacion/0.Envios TS/1 Directorio   (unchanged, just reindexed)
#   46 ToCC                                    (unchanged, just reindexed)
#   47 REPORTE_ESTUDIOC_ESPANIA_04_202
#   48 .../2022/05 Envio Semana 05/ESPAÑA/BD Consolidada
#   49 Key Word
#   50 REPORTE, ESTUDIOC, ESPANIA
#   51 .../2022/05 Envio Semana 05/ESPAÑA/LUXURY
#   52 .../2022/05 Envio Semana 05/ESPAÑA/TRADICIONAL
#   53 Path luxury File when it is divide
#   54 Path for Traditional File when it is divide
#   55 .../2022/05 Envio Semana 05/ESPAÑA/BD Consolidada/Exportadas

# InputSheetName: REPORTE_ESTUDIOC_ESPANIA_03_202 -> ..._04_202
$ws.Range("B3").Value = "REPORTE_ESTUDIOC_ESPANIA_04_202"

# InputnonPartnerFilePath / BlackListTrackingPath: "...Base de Datos" -> "...BD Consolidada"
$ws.Range("B2").Value = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/ESPAÑA/BD Consolidada"
$ws.Range("B13").Value = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/ESPAÑA/BD Consolidada"

# New "Key Word" column (E) with header + search keywords
$ws.Range("E1").Value = "Key Word"
$ws.Range("E2").Value = "REPORTE, ESTUDIOC, ESPANIA"

# BDVIPfilePath now filled in with the LUXURY export path + new description column
$ws.Range("B7").Value = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/ESPAÑA/LUXURY"

# BDTempFilePath now filled in with the TRADICIONAL export path + new description column
$ws.Range("B10").Value = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/ESPAÑA/TRADICIONAL"

$ws.Range("C7").Value = "Path luxury File when it is divide"
$ws.Range("C10").Value = "Path for Traditional File when it is divide"

# SurveyTradFilePath / SurveyVipFilePath: "...Base de Datos/Exportadas" -> "...BD Consolidada/Exportadas"
$ws.Range("B16").Value = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/ESPAÑA/BD Consolidada/Exportadas"
$ws.Range("B19").Value = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/ESPAÑA/BD Consolidada/Exportadas"

# BDTempFileName value cleared out entirely
$ws.Range("B9").Clear()
$ws.Rows.Item(9).AutoFit()

# New column E needs a width matching the other best-fit columns
$ws.Columns.Item(5).ColumnWidth = 27.08

# Selection moved from A31 to B1
$ws.Activate()
$ws.Range("B1").Select()
